# Daily attendance processing - 2025-11-30 21:47:08
# Normalizes the "Recorded By" column (G) so that when the list of
# recorders contains "System" alongside exactly one human/account email,
# "System" is moved to the end: "System, <email>" -> "<email>, System".
# Rows where "System" appears with "backup@backdoor.com", or rows with a
# different number of names (e.g. just one name, or three names), are left
# untouched - matching the source data exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $parts = $val -split ", "

    if ($parts.Count -eq 2 -and $parts[0] -eq "System" -and $parts[1] -ne "backup@backdoor.com") {
        $cell.Value = "$($parts[1]), $($parts[0])"
    }
}
